# Apply the commit's changes to the phone-list sheet:
#  - clear the stray duplicate "Михаил" text left in B4
#  - remove the stray duplicate "Михаил" value in F4 entirely (row 4 now ends at E4)
#  - append a new record (row 5) for Сергей, including the extra G/H/I columns
#  - leave the selection where the editor left off (C10)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 clean-up -------------------------------------------------
# B4 is blanked out but keeps its (now-empty) text entry - a lone "'" is
# Excel's classic "force text, empty content" entry so the cell stays a
# text cell instead of turning into a fully blank/general cell.
$ws.Range("B4").Value = "'"
$ws.Range("F4").Clear()

# --- New row 5: Сергей's record -------------------------------------
# A5 holds a long numeric-looking phone number; force text so it is stored
# the same way the other phone numbers in column A are (shared string, not
# a number).
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "89038452680"

$ws.Range("B5").Value = "Сергей"
$ws.Range("C5").Value = "M"
$ws.Range("D5").Value = "ПАО ""Вымпел-Коммуникации"""
$ws.Range("E5").Value = "Тульская обл."
$ws.Range("F5").Value = "cthutq"
$ws.Range("G5").Value = "vip.klim1964@mail.ru "
$ws.Range("H5").Value = "Россия"
$ws.Range("I5").Value = "Заемщик"

# --- Selection, matching where the author left the cursor -----------
$ws.Range("C10").Select()
